# Updated cryptocurrency price/volume data (symbol list) for Sheet1
# Commit: Updated symbol list on Wed Feb  1 13:46:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'309.24"
$ws.Range('E2').Value = "'-0.55%"
$ws.Range('D3').Value = "'37.36"
$ws.Range('E3').Value = "'-0.38%"
$ws.Range('E4').Value = "'0.39%"
$ws.Range('D5').Value = "'0.07845"
$ws.Range('E5').Value = "'-0.10%"
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D6').Value = "'8.259"
$ws.Range('E6').Value = "'0.53%"
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = "'1.880"
$ws.Range('E7').Value = "'-1.18%"
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = "'2.972"
$ws.Range('E8').Value = "'9.11%"
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = "'0.9230"
$ws.Range('E9').Value = "'-0.41%"
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1110"
$ws.Range('E10').Value = "'-7.38%"
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1908"
$ws.Range('E11').Value = "'0.29%"
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.08884"
$ws.Range('E12').Value = "'-6.07%"
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03337"
$ws.Range('E13').Value = "'-2.33%"
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09607"
$ws.Range('E14').Value = "'-0.11%"
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001385"
$ws.Range('E15').Value = "'1.10%"
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.006001"
$ws.Range('E16').Value = "'2.49%"
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.395"
$ws.Range('E17').Value = "'-4.01%"
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = "'4.392"
$ws.Range('E18').Value = "'0.00%"
$ws.Range('D19').Value = "'0.3456"
$ws.Range('E19').Value = "'0.91%"
$ws.Range('D20').Value = "'6.366"
$ws.Range('E20').Value = "'21.14%"
$ws.Range('D21').Value = "'0.1314"
$ws.Range('E21').Value = "'3.16%"
$ws.Range('E22').Value = "'-6.99%"
$ws.Range('D23').Value = "'0.04350"
$ws.Range('E23').Value = "'0.13%"
$ws.Range('E24').Value = "'0.56%"
$ws.Range('D25').Value = "'0.004284"
$ws.Range('E25').Value = "'0.28%"
$ws.Range('D26').Value = "'0.0001402"
$ws.Range('E26').Value = "'8.12%"
$ws.Range('D27').Value = "'0.0002901"
$ws.Range('D39').Value = "'0.02168"
$ws.Range('E39').Value = "'3.60%"
$ws.Range('D40').Value = "'0.05022"
$ws.Range('E40').Value = "'-0.60%"
$ws.Range('D41').Value = "'0.007583"
$ws.Range('E41').Value = "'-0.59%"
$ws.Range('D42').Value = "'0.1356"
$ws.Range('E42').Value = "'0.45%"
$ws.Range('D43').Value = "'0.008508"
$ws.Range('E43').Value = "'-6.60%"
$ws.Range('D44').Value = "'0.002071"
$ws.Range('E44').Value = "'3.51%"
$ws.Range('D45').Value = "'0.008144"
$ws.Range('E45').Value = "'-5.00%"
$ws.Range('D46').Value = "'0.00006518"
$ws.Range('E46').Value = "'-2.56%"
$ws.Range('E47').Value = "'0.28%"
$ws.Range('D48').Value = "'0.003298"
$ws.Range('E48').Value = "'13.87%"
$ws.Range('E49').Value = "'20.58%"
$ws.Range('D50').Value = "'0.00002101"
$ws.Range('E50').Value = "'0.28%"
$ws.Range('D51').Value = "'0.0002001"
$ws.Range('E51').Value = "'0.28%"
